$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Duplicate the formatting of the last existing data row (row 3) onto the new
# row 4 so that the date cell (G4) picks up the same style (date number format)
# as G3 without introducing a brand new style entry.
$ws.Range("A3:H3").Copy() | Out-Null
$ws.Range("A4:H4").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Populate the new trade row (row 4) with the recorded trade data.
$ws.Range("A4").Value = 10146.5
$ws.Range("B4").Value = 10055
$ws.Range("C4").Value = 107.96
$ws.Range("D4").Value = 108.94
$ws.Range("E4").Value = $false
$ws.Range("F4").Value = 0.91
$ws.Range("G4").Value = 42609.505590277775
$ws.Range("H4").Value = $true
